$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = "G Squad"
$ws.Range("B27").Value = "JJPVG0GJ"

$ws.Range("A28").Value = "Delivering Chaos"
$ws.Range("B28").Value = "YRVGGUG8"

$ws.Range("A29").Value = "g-solo"
$ws.Range("B29").Value = "QR9PU9CU"

$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1

$ws.Range("A30").Select() | Out-Null
